$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.745.34"
$ws.Range("E2").Value = "  -2.26%  "

# Row 3
$ws.Range("D3").Value = "1.798.23"
$ws.Range("E3").Value = "  -1.65%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'308.78"

# Row 6
$ws.Range("E6").Value = "  -0.04%  "

# Row 7
$ws.Range("D7").Value = "'0.4614"
$ws.Range("E7").Value = "  +3.11%  "

# Row 8
$ws.Range("D8").Value = "'0.3716"
$ws.Range("E8").Value = "  -1.33%  "

# Row 9
$ws.Range("D9").Value = "'0.07248"
$ws.Range("E9").Value = "  -3.74%  "

# Row 10
$ws.Range("D10").Value = "'0.8549"
$ws.Range("E10").Value = "  -4.46%  "

# Row 11
$ws.Range("E11").Value = "  -3.30%  "

# Row 12
$ws.Range("D12").Value = "1.766.82"
$ws.Range("E12").Value = "  -3.29%  "

# Row 13
$ws.Range("D13").Value = "'5.307"
$ws.Range("E13").Value = "  -2.15%  "

# Row 14
$ws.Range("D14").Value = "'6.490"
$ws.Range("E14").Value = "  -3.69%  "

# Row 15
$ws.Range("D15").Value = "'0.07041"
$ws.Range("E15").Value = "  -1.21%  "

# Row 16
$ws.Range("D16").Value = "'90.53"
$ws.Range("E16").Value = "  -4.31%  "

# Row 17
$ws.Range("E17").Value = "  -0.04%  "

# Row 18
$ws.Range("D18").Value = "'0.000008629"
$ws.Range("E18").Value = "  -2.31%  "

# Row 19
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -0.05%  "

# Row 20
$ws.Range("D20").Value = "'14.61"
$ws.Range("E20").Value = "  -4.38%  "

# Row 21
$ws.Range("D21").Value = "26.760.29"
$ws.Range("E21").Value = "  -2.29%  "

# Row 22
$ws.Range("D22").Value = "'5.290"
$ws.Range("E22").Value = "  +0.04%  "

# Row 23
$ws.Range("E23").Value = "  -2.87%  "

# Row 24
$ws.Range("D24").Value = "2.007.29"
$ws.Range("E24").Value = "  -1.90%  "

# Row 25
$ws.Range("D25").Value = "'1.907"
$ws.Range("E25").Value = "  -4.81%  "

# Row 26
$ws.Range("D26").Value = "'149.91"
$ws.Range("E26").Value = "  -1.11%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'18.15"
$ws.Range("E27").Value = "  -2.52%  "

# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.141"
$ws.Range("E28").Value = "  -12.52%  "

# Row 29
$ws.Range("E29").Value = "  -3.19%  "

# Row 30
$ws.Range("D30").Value = "'114.11"
$ws.Range("E30").Value = "  -3.21%  "

# Row 31
$ws.Range("D31").Value = "'0.08900"
$ws.Range("E31").Value = "  +0.57%  "

# Row 32
$ws.Range("D32").Value = "'0.7537"
$ws.Range("E32").Value = "  -4.29%  "

# Row 33
$ws.Range("D33").Value = "'1.157"
$ws.Range("E33").Value = "  -3.99%  "

# Row 34
$ws.Range("D34").Value = "'4.428"
$ws.Range("E34").Value = "  -2.67%  "

# Row 35
$ws.Range("D35").Value = "'2.888"
$ws.Range("E35").Value = "  -0.12%  "

# Row 36
$ws.Range("D36").Value = "'1.0000"
$ws.Range("E36").Value = "  -0.05%  "

# Row 37
$ws.Range("D37").Value = "'1.116"
$ws.Range("E37").Value = "  +0.52%  "

# Row 38
$ws.Range("D38").Value = "'0.01938"
$ws.Range("E38").Value = "  -2.74%  "

# Row 39
$ws.Range("D39").Value = "'0.05204"
$ws.Range("E39").Value = "  -2.46%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'2.372"
$ws.Range("E40").Value = "  +2.96%  "

# Row 41
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.897"
$ws.Range("E41").Value = "  +0.82%  "

# Row 42
$ws.Range("D42").Value = "'7.159"
$ws.Range("E42").Value = "  -3.27%  "

# Row 43
$ws.Range("D43").Value = "'0.5221"
$ws.Range("E43").Value = "  -1.95%  "

# Row 44
$ws.Range("D44").Value = "'0.1645"
$ws.Range("E44").Value = "  -5.16%  "

# Row 45
$ws.Range("E45").Value = "  -3.72%  "

# Row 46
$ws.Range("D46").Value = "'0.5005"
$ws.Range("E46").Value = "  -3.26%  "

# Row 47
$ws.Range("D47").Value = "'10.26"
$ws.Range("E47").Value = "  -4.91%  "

# Row 48
$ws.Range("D48").Value = "'104.10"
$ws.Range("E48").Value = "  -2.23%  "

# Row 49
$ws.Range("D49").Value = "'0.9999"
$ws.Range("E49").Value = "  -0.07%  "

# Row 50
$ws.Range("E50").Value = "  -3.93%  "

# Row 51
$ws.Range("D51").Value = "'0.06285"
$ws.Range("E51").Value = "  -1.44%  "
